$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.423.56'
$ws.Range("E2").Value = '  +4.52%  '
$ws.Range("D3").Value = '3.838.74'
$ws.Range("E3").Value = '  +9.63%  '
$ws.Range("E4").Value = '  -0.62%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '426.15'
$ws.Range("E5").Value = '  +9.10%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '131.53'
$ws.Range("E6").Value = '  +8.37%  '
$ws.Range("D7").Value = '3.829.42'
$ws.Range("E7").Value = '  +9.77%  '
$ws.Range("E8").Value = '  +5.19%  '
$ws.Range("E9").Value = '  -0.13%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.732'
$ws.Range("E10").Value = '  +9.46%  '
$ws.Range("E11").Value = '  +4.88%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000340'
$ws.Range("E12").Value = '  +1.45%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '41.92'
$ws.Range("E13").Value = '  +8.96%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '10.53'
$ws.Range("E14").Value = '  +15.31%  '
$ws.Range("D15").Value = '4.452.98'
$ws.Range("E15").Value = '  +9.73%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '15.81'
$ws.Range("E16").Value = '  +27.39%  '
$ws.Range("B17").Value = 'TRON'
$ws.Range("C17").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.138'
$ws.Range("E17").Value = '  +1.34%  '
$ws.Range("B18").Value = 'WrappedEther'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D18").Value = '3.766.73'
$ws.Range("E18").Value = '  +7.76%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '20.10'
$ws.Range("E19").Value = '  +8.31%  '
$ws.Range("E20").Value = '  +8.83%  '
$ws.Range("D21").Value = '66.706.78'
$ws.Range("E21").Value = '  +4.84%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '415.61'
$ws.Range("E22").Value = '  +6.42%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '15.10'
$ws.Range("E23").Value = '  +9.84%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '84.92'
$ws.Range("E24").Value = '  +6.04%  '
$ws.Range("E25").Value = '  +8.77%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '37.58'
$ws.Range("E26").Value = '  +14.38%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.02'
$ws.Range("E27").Value = '  +15.04%  '
$ws.Range("E28").Value = '  +9.93%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.34'
$ws.Range("E29").Value = '  +1.86%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.18'
$ws.Range("E30").Value = '  +35.77%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '718.92'
$ws.Range("E31").Value = '  +10.03%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '13.89'
$ws.Range("E32").Value = '  +17.01%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.125'
$ws.Range("E33").Value = '  +15.09%  '
$ws.Range("E34").Value = '  +6.35%  '
$ws.Range("B35").Value = 'NEARProtocol'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.83'
$ws.Range("E35").Value = '  +46.05%  '
$ws.Range("B36").Value = 'Dai'
$ws.Range("C36").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.999'
$ws.Range("E36").Value = '  -0.08%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '39.21'
$ws.Range("E37").Value = '  +7.73%  '
$ws.Range("E38").Value = '  +1.07%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '55.69'
$ws.Range("E39").Value = '  +3.66%  '
$ws.Range("D40").Value = '0.0₃0748'
$ws.Range("E40").Value = '  +18.56%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0465'
$ws.Range("E41").Value = '  +7.17%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.91'
$ws.Range("E42").Value = '  +9.00%  '
$ws.Range("E43").Value = '  +0.29%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.25'
$ws.Range("E44").Value = '  +7.36%  '
$ws.Range("E45").Value = '  +4.94%  '
$ws.Range("E46").Value = '  +10.96%  '
$ws.Range("E47").Value = '  +15.69%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '141.92'
$ws.Range("E48").Value = '  +1.78%  '
$ws.Range("B49").Value = 'Stacks'
$ws.Range("C49").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.84'
$ws.Range("E49").Value = '  +6.12%  '
$ws.Range("B50").Value = 'ARBITRUM'
$ws.Range("C50").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.05'
$ws.Range("E50").Value = '  +5.87%  '
$ws.Range("E51").Value = '  +4.07%  '
